# Atualização automática de BUTIA.xlsx
# - Renomeia "Paineis DARQ" para "PAINEIS DARQ"
# - Renomeia "Recolhimento x Eliminacao" para "RECOLHIMENTO X ELIMINAÇÃO"
# - Remove a planilha "Desarquivamentos Pendentes"

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true
